$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 41.428665
$ws.Range("H2").Value = 124.285995
$ws.Range("I2").Value = 0.06969137269740189
$ws.Range("J2").Value = 0.06969137269740189
$ws.Range("M2").Value = 159.4836373333333
$ws.Range("N2").Value = 478.450912
$ws.Range("O2").Value = 0.2983285084902258
$ws.Range("P2").Value = 0.2983285084902258
$ws.Range("Q2").Value = 6607.19418406416
$ws.Range("R2").Value = 59464.74765657743
$ws.Range("S2").Value = 0.02079092327145235
$ws.Range("T2").Value = 0.02079092327145235

$ws.Range("G3").Value = 41.428665
$ws.Range("H3").Value = 124.285995
$ws.Range("I3").Value = 0.06969137269740189
$ws.Range("J3").Value = 0.06969137269740189
$ws.Range("O3").Value = 0.3227862111630279
$ws.Range("P3").Value = 0.3227862111630279
$ws.Range("Q3").Value = 7148.868165116503
$ws.Range("R3").Value = 64339.81348604853
$ws.Range("S3").Value = 0.02249541414374484
$ws.Range("T3").Value = 0.02249541414374484

$ws.Range("G4").Value = 41.428665
$ws.Range("H4").Value = 124.285995
$ws.Range("I4").Value = 0.06969137269740189
$ws.Range("J4").Value = 0.06969137269740189
$ws.Range("M4").Value = 74.38770566666666
$ws.Range("N4").Value = 223.163117
$ws.Range("O4").Value = 0.1391489036280481
$ws.Range("P4").Value = 0.1391489036280482
$ws.Range("Q4").Value = 3081.783338182935
$ws.Range("R4").Value = 27736.05004364641
$ws.Range("S4").Value = 0.00969747810317716
$ws.Range("T4").Value = 0.009697478103177162

$ws.Range("G5").Value = 41.428665
$ws.Range("H5").Value = 124.285995
$ws.Range("I5").Value = 0.06969137269740189
$ws.Range("J5").Value = 0.06969137269740189
$ws.Range("M5").Value = 58.41461433333333
$ws.Range("N5").Value = 175.243843
$ws.Range("O5").Value = 0.1092697975759847
$ws.Range("P5").Value = 0.1092697975759848
$ws.Range("Q5").Value = 2420.039488319865
$ws.Range("R5").Value = 21780.35539487878
$ws.Range("S5").Value = 0.007615162187437614
$ws.Range("T5").Value = 0.007615162187437616

$ws.Range("G6").Value = 41.428665
$ws.Range("H6").Value = 124.285995
$ws.Range("I6").Value = 0.06969137269740189
$ws.Range("J6").Value = 0.06969137269740189
$ws.Range("M6").Value = 69.746216
$ws.Range("N6").Value = 209.238648
$ws.Range("O6").Value = 0.1304665791427133
$ws.Range("P6").Value = 0.1304665791427133
$ws.Range("Q6").Value = 2889.49261768164
$ws.Range("R6").Value = 26005.43355913476
$ws.Range("S6").Value = 0.009092394991589913
$ws.Range("T6").Value = 0.009092394991589915

$ws.Range("H7").Value = 510.696747
$ws.Range("I7").Value = 0.2863649869040173
$ws.Range("J7").Value = 0.2863649869040173
$ws.Range("M7").Value = 159.4836373333333
$ws.Range("N7").Value = 478.450912
$ws.Range("O7").Value = 0.2983285084902258
$ws.Range("P7").Value = 0.2983285084902258
$ws.Range("Q7").Value = 27149.2582619537
$ws.Range("R7").Value = 244343.3243575832
$ws.Range("S7").Value = 0.08543083942689852
$ws.Range("T7").Value = 0.08543083942689854

$ws.Range("H8").Value = 510.696747
$ws.Range("I8").Value = 0.2863649869040173
$ws.Range("J8").Value = 0.2863649869040173
$ws.Range("O8").Value = 0.3227862111630279
$ws.Range("P8").Value = 0.3227862111630279
$ws.Range("Q8").Value = 29375.02102836975
$ws.Range("R8").Value = 264375.1892553277
$ws.Range("S8").Value = 0.09243466913249784
$ws.Range("T8").Value = 0.09243466913249787

$ws.Range("H9").Value = 510.696747
$ws.Range("I9").Value = 0.2863649869040173
$ws.Range("J9").Value = 0.2863649869040173
$ws.Range("M9").Value = 74.38770566666666
$ws.Range("N9").Value = 223.163117
$ws.Range("O9").Value = 0.1391489036280481
$ws.Range("P9").Value = 0.1391489036280482
$ws.Range("Q9").Value = 12663.18643358671
$ws.Range("R9").Value = 113968.6779022804
$ws.Range("S9").Value = 0.03984737396515436
$ws.Range("T9").Value = 0.03984737396515439

$ws.Range("H10").Value = 510.696747
$ws.Range("I10").Value = 0.2863649869040173
$ws.Range("J10").Value = 0.2863649869040173
$ws.Range("M10").Value = 58.41461433333333
$ws.Range("N10").Value = 175.243843
$ws.Range("O10").Value = 0.1092697975759847
$ws.Range("P10").Value = 0.1092697975759848
$ws.Range("Q10").Value = 9944.051172430969
$ws.Range("R10").Value = 89496.46055187871
$ws.Range("S10").Value = 0.03129104415185149
$ws.Range("T10").Value = 0.0312910441518515

$ws.Range("H11").Value = 510.696747
$ws.Range("I11").Value = 0.2863649869040173
$ws.Range("J11").Value = 0.2863649869040173
$ws.Range("M11").Value = 69.746216
$ws.Range("N11").Value = 209.238648
$ws.Range("O11").Value = 0.1304665791427133
$ws.Range("P11").Value = 0.1304665791427133
$ws.Range("Q11").Value = 11873.05520891978
$ws.Range("R11").Value = 106857.496880278
$ws.Range("S11").Value = 0.03736106022761503
$ws.Range("T11").Value = 0.03736106022761505

$ws.Range("G12").Value = 244.5761666666666
$ws.Range("H12").Value = 733.7284999999999
$ws.Range("I12").Value = 0.4114264551867299
$ws.Range("J12").Value = 0.41142645518673
$ws.Range("M12").Value = 159.4836373333333
$ws.Range("N12").Value = 478.450912
$ws.Range("O12").Value = 0.2983285084902258
$ws.Range("P12").Value = 0.2983285084902258
$ws.Range("Q12").Value = 39005.89666504355
$ws.Range("R12").Value = 351053.069985392
$ws.Range("S12").Value = 0.1227402407292779
$ws.Range("T12").Value = 0.1227402407292779

$ws.Range("G13").Value = 244.5761666666666
$ws.Range("H13").Value = 733.7284999999999
$ws.Range("I13").Value = 0.4114264551867299
$ws.Range("J13").Value = 0.41142645518673
$ws.Range("O13").Value = 0.3227862111630279
$ws.Range("P13").Value = 0.3227862111630279
$ws.Range("Q13").Value = 42203.69572202149
$ws.Range("R13").Value = 379833.2614981934
$ws.Range("S13").Value = 0.1328027866419599
$ws.Range("T13").Value = 0.1328027866419599

$ws.Range("G14").Value = 244.5761666666666
$ws.Range("H14").Value = 733.7284999999999
$ws.Range("I14").Value = 0.4114264551867299
$ws.Range("J14").Value = 0.41142645518673
$ws.Range("M14").Value = 74.38770566666666
$ws.Range("N14").Value = 223.163117
$ws.Range("O14").Value = 0.1391489036280481
$ws.Range("P14").Value = 0.1391489036280482
$ws.Range("Q14").Value = 18193.45989908161
$ws.Range("R14").Value = 163741.1390917345
$ws.Range("S14").Value = 0.05724954016280775
$ws.Range("T14").Value = 0.05724954016280777

$ws.Range("G15").Value = 244.5761666666666
$ws.Range("H15").Value = 733.7284999999999
$ws.Range("I15").Value = 0.4114264551867299
$ws.Range("J15").Value = 0.41142645518673
$ws.Range("M15").Value = 58.41461433333333
$ws.Range("N15").Value = 175.243843
$ws.Range("O15").Value = 0.1092697975759847
$ws.Range("P15").Value = 0.1092697975759848
$ws.Range("Q15").Value = 14286.82245095839
$ws.Range("R15").Value = 128581.4020586255
$ws.Range("S15").Value = 0.04495648547565894
$ws.Range("T15").Value = 0.04495648547565895

$ws.Range("G16").Value = 244.5761666666666
$ws.Range("H16").Value = 733.7284999999999
$ws.Range("I16").Value = 0.4114264551867299
$ws.Range("J16").Value = 0.41142645518673
$ws.Range("M16").Value = 69.746216
$ws.Range("N16").Value = 209.238648
$ws.Range("O16").Value = 0.1304665791427133
$ws.Range("P16").Value = 0.1304665791427133
$ws.Range("Q16").Value = 17058.26214878533
$ws.Range("R16").Value = 153524.359339068
$ws.Range("S16").Value = 0.05367740217702549
$ws.Range("T16").Value = 0.05367740217702551

$ws.Range("G17").Value = 24.173247
$ws.Range("H17").Value = 72.51974100000001
$ws.Range("I17").Value = 0.04066427836821081
$ws.Range("J17").Value = 0.04066427836821081
$ws.Range("M17").Value = 159.4836373333333
$ws.Range("N17").Value = 478.450912
$ws.Range("O17").Value = 0.2983285084902258
$ws.Range("P17").Value = 0.2983285084902258
$ws.Range("Q17").Value = 3855.237357717089
$ws.Range("R17").Value = 34697.1362194538
$ws.Range("S17").Value = 0.01213131351441969
$ws.Range("T17").Value = 0.01213131351441969

$ws.Range("G18").Value = 24.173247
$ws.Range("H18").Value = 72.51974100000001
$ws.Range("I18").Value = 0.04066427836821081
$ws.Range("J18").Value = 0.04066427836821081
$ws.Range("O18").Value = 0.3227862111630279
$ws.Range("P18").Value = 0.3227862111630279
$ws.Range("Q18").Value = 4171.299169929759
$ws.Range("R18").Value = 37541.69252936784
$ws.Range("S18").Value = 0.01312586834415344
$ws.Range("T18").Value = 0.01312586834415344

$ws.Range("G19").Value = 24.173247
$ws.Range("H19").Value = 72.51974100000001
$ws.Range("I19").Value = 0.04066427836821081
$ws.Range("J19").Value = 0.04066427836821081
$ws.Range("M19").Value = 74.38770566666666
$ws.Range("N19").Value = 223.163117
$ws.Range("O19").Value = 0.1391489036280481
$ws.Range("P19").Value = 0.1391489036280482
$ws.Range("Q19").Value = 1798.192382843633
$ws.Range("R19").Value = 16183.7314455927
$ws.Range("S19").Value = 0.005658389751762289
$ws.Range("T19").Value = 0.00565838975176229

$ws.Range("G20").Value = 24.173247
$ws.Range("H20").Value = 72.51974100000001
$ws.Range("I20").Value = 0.04066427836821081
$ws.Range("J20").Value = 0.04066427836821081
$ws.Range("M20").Value = 58.41461433333333
$ws.Range("N20").Value = 175.243843
$ws.Range("O20").Value = 0.1092697975759847
$ws.Range("P20").Value = 0.1092697975759848
$ws.Range("Q20").Value = 1412.070900689407
$ws.Range("R20").Value = 12708.63810620466
$ws.Range("S20").Value = 0.004443377465867891
$ws.Range("T20").Value = 0.004443377465867891

$ws.Range("G21").Value = 24.173247
$ws.Range("H21").Value = 72.51974100000001
$ws.Range("I21").Value = 0.04066427836821081
$ws.Range("J21").Value = 0.04066427836821081
$ws.Range("M21").Value = 69.746216
$ws.Range("N21").Value = 209.238648
$ws.Range("O21").Value = 0.1304665791427133
$ws.Range("P21").Value = 0.1304665791427133
$ws.Range("Q21").Value = 1685.992506683352
$ws.Range("R21").Value = 15173.93256015017
$ws.Range("S21").Value = 0.005305329292007501
$ws.Range("T21").Value = 0.005305329292007503

$ws.Range("G22").Value = 114.0486906666667
$ws.Range("H22").Value = 342.146072
$ws.Range("I22").Value = 0.19185290684364
$ws.Range("J22").Value = 0.19185290684364
$ws.Range("M22").Value = 159.4836373333333
$ws.Range("N22").Value = 478.450912
$ws.Range("O22").Value = 0.2983285084902258
$ws.Range("P22").Value = 0.2983285084902258
$ws.Range("Q22").Value = 18188.90002062419
$ws.Range("R22").Value = 163700.1001856177
$ws.Range("S22").Value = 0.05723519154817735
$ws.Range("T22").Value = 0.05723519154817735

$ws.Range("G23").Value = 114.0486906666667
$ws.Range("H23").Value = 342.146072
$ws.Range("I23").Value = 0.19185290684364
$ws.Range("J23").Value = 0.19185290684364
$ws.Range("O23").Value = 0.3227862111630279
$ws.Range("P23").Value = 0.3227862111630279
$ws.Range("Q23").Value = 19680.07064625793
$ws.Range("R23").Value = 177120.6358163213
$ws.Range("S23").Value = 0.0619274729006719
$ws.Range("T23").Value = 0.0619274729006719

$ws.Range("G24").Value = 114.0486906666667
$ws.Range("H24").Value = 342.146072
$ws.Range("I24").Value = 0.19185290684364
$ws.Range("J24").Value = 0.19185290684364
$ws.Range("M24").Value = 74.38770566666666
$ws.Range("N24").Value = 223.163117
$ws.Range("O24").Value = 0.1391489036280481
$ws.Range("P24").Value = 0.1391489036280482
$ws.Range("Q24").Value = 8483.820432980714
$ws.Range("R24").Value = 76354.38389682643
$ws.Range("S24").Value = 0.02669612164514655
$ws.Range("T24").Value = 0.02669612164514656

$ws.Range("G25").Value = 114.0486906666667
$ws.Range("H25").Value = 342.146072
$ws.Range("I25").Value = 0.19185290684364
$ws.Range("J25").Value = 0.19185290684364
$ws.Range("M25").Value = 58.41461433333333
$ws.Range("N25").Value = 175.243843
$ws.Range("O25").Value = 0.1092697975759847
$ws.Range("P25").Value = 0.1092697975759848
$ws.Range("Q25").Value = 6662.110280514967
$ws.Range("R25").Value = 59958.9925246347
$ws.Range("S25").Value = 0.0209637282951688
$ws.Range("T25").Value = 0.0209637282951688

$ws.Range("G26").Value = 114.0486906666667
$ws.Range("H26").Value = 342.146072
$ws.Range("I26").Value = 0.19185290684364
$ws.Range("J26").Value = 0.19185290684364
$ws.Range("M26").Value = 69.746216
$ws.Range("N26").Value = 209.238648
$ws.Range("O26").Value = 0.1304665791427133
$ws.Range("P26").Value = 0.1304665791427133
$ws.Range("Q26").Value = 7954.464613754518
$ws.Range("R26").Value = 71590.18152379066
$ws.Range("S26").Value = 0.02503039245447536
$ws.Range("T26").Value = 0.02503039245447537
